# Update column AF (average_doctor_MA) values for rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.726
$ws.Range("AF5").Value = 0.951
$ws.Range("AF6").Value = 0.823
$ws.Range("AF7").Value = 0.895
$ws.Range("AF8").Value = 0.867
$ws.Range("AF9").Value = 0.732
$ws.Range("AF10").Value = 0.951
$ws.Range("AF11").Value = 0.951
$ws.Range("AF12").Value = 1.256
$ws.Range("AF13").Value = 1.634
